# Auto-generated COM-interop script updating Leve profit-analysis values
# across the ALC/ARM/BSM/CRP/CUL/LTW/WVR sheets (scheduled price-refresh run).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4341.4287
$ws.Range("I74").Value = 3955.5557
$ws.Range("K74").Value = 3955.5557
$ws.Range("M74").Value = -3019.5557

$ws.Range("H77").Value = 4341.4287
$ws.Range("I77").Value = 3955.5557
$ws.Range("K77").Value = 19777.7785
$ws.Range("M77").Value = -15097.7785

$ws.Range("H94").Value = 2750
$ws.Range("I94").Value = 2750
$ws.Range("K94").Value = 2750
$ws.Range("M94").Value = -2299

$ws.Range("H96").Value = 493.66666
$ws.Range("I96").Value = 226
$ws.Range("J96").Value = 1029
$ws.Range("K96").Value = 678
$ws.Range("L96").Value = 3087
$ws.Range("M96").Value = 695
$ws.Range("N96").Value = -5833

$ws.Range("H99").Value = 797.4286
$ws.Range("I99").Value = 347
$ws.Range("J99").Value = 3500
$ws.Range("K99").Value = 1041
$ws.Range("L99").Value = 10500
$ws.Range("M99").Value = 457
$ws.Range("N99").Value = -13496

$ws.Range("H106").Value = 5254.6924
$ws.Range("I106").Value = 8222
$ws.Range("J106").Value = 3400.125
$ws.Range("K106").Value = 8222
$ws.Range("L106").Value = 3400.125
$ws.Range("M106").Value = -7591
$ws.Range("N106").Value = -4662.125

$ws.Range("H132").Value = 8327.947
$ws.Range("I132").Value = 8424.826
$ws.Range("K132").Value = 25274.478
$ws.Range("M132").Value = -22744.478

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8835.304
$ws.Range("I32").Value = 7784.7573
$ws.Range("J32").Value = 17006.223
$ws.Range("K32").Value = 7784.7573
$ws.Range("L32").Value = 17006.223
$ws.Range("M32").Value = -7497.7573
$ws.Range("N32").Value = -17580.223

$ws.Range("H132").Value = 704095.25
$ws.Range("J132").Value = 5558.645
$ws.Range("L132").Value = 16675.935
$ws.Range("N132").Value = -21735.935

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 892.0357
$ws.Range("I94").Value = 811.7619
$ws.Range("J94").Value = 1132.8572
$ws.Range("K94").Value = 811.7619
$ws.Range("L94").Value = 1132.8572
$ws.Range("M94").Value = -360.7619
$ws.Range("N94").Value = -2034.8572

$ws.Range("H124").Value = 45495
$ws.Range("J124").Value = 45495
$ws.Range("L124").Value = 45495
$ws.Range("N124").Value = -55315

$ws.Range("H134").Value = 6391.5312
$ws.Range("I134").Value = 2414.45
$ws.Range("J134").Value = 13020
$ws.Range("K134").Value = 7243.349999999999
$ws.Range("L134").Value = 39060
$ws.Range("M134").Value = -4708.349999999999
$ws.Range("N134").Value = -44130

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6805311
$ws.Range("I31").Value = 1952
$ws.Range("J31").Value = 14496065
$ws.Range("K31").Value = 1952
$ws.Range("L31").Value = 14496065
$ws.Range("M31").Value = -1657
$ws.Range("N31").Value = -14496655

$ws.Range("H34").Value = 6805311
$ws.Range("I34").Value = 1952
$ws.Range("J34").Value = 14496065
$ws.Range("K34").Value = 1952
$ws.Range("L34").Value = 14496065
$ws.Range("M34").Value = -1750
$ws.Range("N34").Value = -14496469

$ws.Range("H105").Value = 1447.95
$ws.Range("I105").Value = 1447.95
$ws.Range("K105").Value = 1447.95
$ws.Range("M105").Value = 299.05

$ws.Range("H134").Value = 1218.7727
$ws.Range("I134").Value = 743.8421
$ws.Range("J134").Value = 4226.6665
$ws.Range("K134").Value = 2231.5263
$ws.Range("L134").Value = 12679.9995
$ws.Range("M134").Value = 303.4737
$ws.Range("N134").Value = -17749.9995

$ws.Range("H135").Value = 38000
$ws.Range("J135").Value = 38000
$ws.Range("L135").Value = 38000
$ws.Range("N135").Value = -48140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 478.5909
$ws.Range("I5").Value = 214.94444
$ws.Range("K5").Value = 644.83332
$ws.Range("M5").Value = -532.83332

$ws.Range("H132").Value = 1152.2916
$ws.Range("I132").Value = 975.8823
$ws.Range("J132").Value = 1580.7142
$ws.Range("K132").Value = 8782.9407
$ws.Range("L132").Value = 14226.4278
$ws.Range("M132").Value = -6252.940699999999
$ws.Range("N132").Value = -19286.4278

$ws.Range("H135").Value = 478.5909
$ws.Range("I135").Value = 214.94444
$ws.Range("K135").Value = 1934.49996
$ws.Range("M135").Value = 600.5000400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1729.1428
$ws.Range("I7").Value = 1767.3334
$ws.Range("K7").Value = 1767.3334
$ws.Range("M7").Value = -1655.3334

$ws.Range("H55").Value = 343.05884
$ws.Range("I55").Value = 448.8889
$ws.Range("J55").Value = 224
$ws.Range("K55").Value = 448.8889
$ws.Range("L55").Value = 224
$ws.Range("M55").Value = -275.8889
$ws.Range("N55").Value = -570

$ws.Range("H93").Value = 1803.44
$ws.Range("I93").Value = 1530.6316
$ws.Range("J93").Value = 2667.3333
$ws.Range("K93").Value = 1530.6316
$ws.Range("L93").Value = 2667.3333
$ws.Range("M93").Value = -282.6315999999999
$ws.Range("N93").Value = -5163.3333

$ws.Range("H126").Value = 1729.1428
$ws.Range("I126").Value = 1767.3334
$ws.Range("K126").Value = 5302.0002
$ws.Range("M126").Value = -2832.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5134
$ws.Range("I62").Value = 3416
$ws.Range("K62").Value = 3416
$ws.Range("M62").Value = -2792

$ws.Range("H65").Value = 5134
$ws.Range("I65").Value = 3416
$ws.Range("K65").Value = 17080
$ws.Range("M65").Value = -13960

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H81").Value = 2562.75
$ws.Range("I81").Value = 2917
$ws.Range("J81").Value = 1500
$ws.Range("K81").Value = 5834
$ws.Range("L81").Value = 3000
$ws.Range("M81").Value = -4773
$ws.Range("N81").Value = -5122

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H84").Value = 2562.75
$ws.Range("I84").Value = 2917
$ws.Range("J84").Value = 1500
$ws.Range("K84").Value = 29170
$ws.Range("L84").Value = 15000
$ws.Range("M84").Value = -23866
$ws.Range("N84").Value = -25608

$ws.Range("H126").Value = 2532.5173
$ws.Range("I126").Value = 3047.389
$ws.Range("J126").Value = 1690
$ws.Range("K126").Value = 9142.167000000001
$ws.Range("L126").Value = 5070
$ws.Range("M126").Value = -6672.167000000001
$ws.Range("N126").Value = -10010
